$wb = $excel.ActiveWorkbook

# --- Sheet: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Penalty/Reward adjustment: reduce requested quantity for the week of 2023-06-11
$ws1.Range("B11").Value = 50

# Remove the last two rows of data (weeks of 2024-03-10 and 2024-03-17)
$ws1.Rows.Item(21).Delete()
$ws1.Rows.Item(20).Delete()

# --- Sheet: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Penalty/Reward adjustment: reduce requested quantity for the month of 2023-06
$ws2.Range("B5").Value = 50

# Remove the last row of data (month of 2024-03-31)
$ws2.Rows.Item(10).Delete()
